# Applies the "Updated cryptos list" data refresh to Sheet1.
# Only the cells whose text actually changed (per the diff) are touched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.718.84"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "3.089.35"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.086.47"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("E13").Value = "  -4.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.55%  "
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "3.604.65"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "66.635.35"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").Value = "3.088.73"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "483.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.687"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.04%  "
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("E27").Value = "  -3.83%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("D34").Value = "0.0₃0930"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.80%  "
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "48.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.309"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.122"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("E42").Value = "  -5.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Value = "2.775.93"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "367.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("E51").Value = "  -2.85%  "
